# Insert a new record row at row 645, shifting existing rows 645-734 down to 646-735.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("645:645").Insert()

# Populate the newly inserted row 645 with the new record's values.
$ws.Range("A645").Value = 3
$ws.Range("B645").Value = "Femacal de La Calera"
$ws.Range("C645").Value = "Coquimbo"
$ws.Range("D645").Value = 45127
$ws.Range("E645").Value = 5
$ws.Range("F645").Value = 100112021
$ws.Range("G645").Value = "Ají"
$ws.Range("H645").Value = "Inferno"
$ws.Range("I645").Value = "Primera"
$ws.Range("J645").Value = 70
$ws.Range("K645").Value = 12000
$ws.Range("L645").Value = 12500
$ws.Range("M645").Value = 12250
$ws.Range("N645").Value = "$/caja 10 kilos"
$ws.Range("O645").Value = "Región de Arica y Parinacota"
$ws.Range("P645").Value = 1225
$ws.Range("Q645").Value = 10
$ws.Range("R645").Value = "Hortaliza"

Write-Output "Row 645 inserted and populated"
